# Add "PersonStateFingerprintIdentification" row to the Post Consolidation
# Identifiers section of the SSP mapping sheet.
#
# The new row is inserted right after the existing "Person State
# Identification ID" row (row 8) and before the "Person FBI Identification
# ID" row (old row 9, which shifts down to row 10), mirroring the same
# three-column layout (Element Name / Element Description / NEIM 3.0
# Mapping) used by every other data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; Excel shifts row 9 (and everything below
# it) down to row 10 and copies the formatting of the row above (row 8) into
# the freshly inserted row, which matches the target style.
$ws.Rows("9:9").Insert()

$ws.Range("A9").Value = "Person State Fingerprint ID"
$ws.Range("B9").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C9").Value = "/CHcr-doc:CriminalHistoryConsolidationReport/nc:Person/CHcr-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"
